$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 2 ---
$ws.Range("D2").Value = 45250
$ws.Range("M2").Value = 30

# --- Update existing row 3 ---
$ws.Range("D3").Value = 45239
$ws.Range("M3").Value = 25
$ws.Range("R3").Value = "Provincia de San Felipe de Aconcagua"

# --- Update existing row 4 ---
$ws.Range("D4").Value = 45244
$ws.Range("M4").Value = 70

# --- Update existing row 5 ---
$ws.Range("D5").Value = 45243
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 35000
$ws.Range("O5").Value = 35000
$ws.Range("P5").Value = 35000
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 7000

# --- Add new row 6 (carries the data that used to be in row 5) ---
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 44179
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101006
$ws.Range("J6").Value = "Breva"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("Q6").Value = "$/bandeja 5 kilos"
$ws.Range("R6").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S6").Value = 3000
$ws.Range("T6").Value = 5

# Match the date-style formatting used by the other rows in column D
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
